$wb = $excel.ActiveWorkbook

# Metadata sheet updates
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# Elements sheet updates
$wsElements = $wb.Worksheets.Item("Elements")
# Row 2 = EIVL_TS -> Definition (column M)
$wsElements.Range("M2").Value = "A quantity specifying a point on the axis of natural time. A point in time is most often represented as a calendar expression."
# Row 5 = EIVL_TS.operator -> Binding Value Set (column Z)
$wsElements.Range("Z5").Value = "http://hl7.org/cda/stds/core/ValueSet/CDASetOperator"
$wsElements.Columns.Item("Z").AutoFit()
